$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The diff inserts two brand-new data rows right after the current row 180
# (SpecCode "SQUIMAN"), pushing every row that used to be 181-241 down to
# 183-241... i.e. down by 2 (they end up at 183-243).
$ws.Rows("181:182").Insert()

# After the insert, the row that used to be 181 (SOLEMON2025/ITA17/28/2-RAP/
# GOBISUE/...) is now row 183. Copy its formatting + text-typed cells (so
# the textual station code "28" etc. stay text, matching the rest of the
# sheet) into the two freshly inserted blank rows, then overwrite just the
# handful of cells that differ for the new SQUIMAN records.
$ws.Rows("183").Copy()
$ws.Rows("181").PasteSpecial(-4104)
$ws.Rows("182").PasteSpecial(-4104)
$excel.CutCopyMode = 0

# Row 181 - new SQUIMAN record (SampN=1, SpecN=1, L=14mm, W=4g, Sex=I)
$ws.Cells.Item(181, 5).Value2 = "SQUIMAN"
$ws.Cells.Item(181, 8).Value2 = 14
$ws.Cells.Item(181, 9).Value2 = 4
$ws.Cells.Item(181, 16).Value2 = "NA l inferred"

# Row 182 - new SQUIMAN record (SampN=1, SpecN=1, L=9mm, W=1g, Sex=I)
$ws.Cells.Item(182, 5).Value2 = "SQUIMAN"
$ws.Cells.Item(182, 8).Value2 = 9
$ws.Cells.Item(182, 9).Value2 = 1
$ws.Cells.Item(182, 16).Value2 = "NA l inferred"
